$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp (2024-06-13 09:46:50)
$ws.Name = "IClientBalance-20240613-094650-"

# The whole export was refreshed a day later: every "Dt. Referencia" (column G)
# moves from 45455 (2024-06-12) to 45456 (2024-06-13)
$ws.Range("G2:G257").Value = 45456

# A handful of balances were corrected in the refreshed export
# Row 107: Saldo Previsto / Vl. Total both increase by 10000
$ws.Range("D107").Value = 10056.27
$ws.Range("H107").Value = 10056.27

# Row 146: Saldo Previsto / Vl. Total reset to 0
$ws.Range("D146").Value = 0
$ws.Range("H146").Value = 0

# Row 168: Vl. Projetado cleared to 0, Saldo Previsto reduced to match Vl. Total
$ws.Range("D168").Value = 9389.98
$ws.Range("E168").Value = 0

# Row 191: Saldo Previsto / Vl. Total both increase slightly
$ws.Range("D191").Value = 59.76
$ws.Range("H191").Value = 59.76

# Row 245: Saldo Previsto / Vl. Total both reduced by 3000
$ws.Range("D245").Value = 978.94
$ws.Range("H245").Value = 978.94
